# Add a new activity-log entry (row 2) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column (A2) -- set the number format BEFORE the value so the
# engine doesn't also register an extra (unused) auto date format.
# "mm-dd-yy" resolves to Excel's built-in numFmtId 14 (m/d/yyyy).
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = (Get-Date -Year 2020 -Month 2 -Day 3 -Hour 0 -Minute 0 -Second 0).Date

# UBIT / DURATION / ACTIVITY & TOOL DESCRIPTION columns
$ws.Range("B2").Value = "jaclemon"
$ws.Range("C2").Value = "60 minutes"
$ws.Range("D2").Value = "Setting up Github repo for project"

# Leave the selection where the author left it when they saved.
$ws.Range("D5").Select() | Out-Null
